$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the values (A3:F32) up into (A2:F31) as static values (no formulas),
# then delete the now-duplicate last row (row 32).
$vals = $ws.Range("A3:F32").Value2

$dest = $ws.Range("A2:F31")
$dest.Value2 = $vals

# Remove the old trailing row which is now a duplicate of row 31.
$ws.Rows.Item(32).Delete()

# Update the selection to match the post-edit state.
$ws.Range("G16").Select()
